# Update LR-pairs data with refreshed TPM values (Ccl28-Ccr10 edge table)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Resolving-Mac" sending-cluster rows (rows 8-10) entirely
$ws.Rows("8:10").Delete()

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1395456666666667
$ws.Range("H2").Value = 0.418637
$ws.Range("I2").Value = 0.5708284189068497
$ws.Range("J2").Value = 0.5708284189068498
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.058783666666667
$ws.Range("N2").Value = 3.176351
$ws.Range("O2").Value = 0.5014862149947701
$ws.Range("P2").Value = 0.5014862149947702
$ws.Range("Q2").Value = 0.1477486726207778
$ws.Range("R2").Value = 1.329738053587
$ws.Range("S2").Value = 0.2862625832090451
$ws.Range("T2").Value = 0.2862625832090452

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1395456666666667
$ws.Range("H3").Value = 0.418637
$ws.Range("I3").Value = 0.5708284189068497
$ws.Range("J3").Value = 0.5708284189068498
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.003554
$ws.Range("N3").Value = 3.010662
$ws.Range("O3").Value = 0.475327031240749
$ws.Range("P3").Value = 0.4753270312407492
$ws.Range("Q3").Value = 0.140041611966
$ws.Range("R3").Value = 1.260374507694
$ws.Range("S3").Value = 0.2713301777068435
$ws.Range("T3").Value = 0.2713301777068436

# Row 4: FAPs -> Resolving-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl28"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1395456666666667
$ws.Range("H4").Value = 0.418637
$ws.Range("I4").Value = 0.5708284189068497
$ws.Range("J4").Value = 0.5708284189068498
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.048954
$ws.Range("N4").Value = 0.146862
$ws.Range("O4").Value = 0.02318675376448066
$ws.Range("P4").Value = 0.02318675376448067
$ws.Range("Q4").Value = 0.006831318565999999
$ws.Range("R4").Value = 0.06148186709399999
$ws.Range("S4").Value = 0.01323565799096094
$ws.Range("T4").Value = 0.01323565799096095

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl28"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.104916
$ws.Range("H5").Value = 0.314748
$ws.Range("I5").Value = 0.4291715810931503
$ws.Range("J5").Value = 0.4291715810931503
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.058783666666667
$ws.Range("N5").Value = 3.176351
$ws.Range("O5").Value = 0.5014862149947701
$ws.Range("P5").Value = 0.5014862149947702
$ws.Range("Q5").Value = 0.111083347172
$ws.Range("R5").Value = 0.9997501245480002
$ws.Range("S5").Value = 0.215223631785725
$ws.Range("T5").Value = 0.215223631785725

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ccl28"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.104916
$ws.Range("H6").Value = 0.314748
$ws.Range("I6").Value = 0.4291715810931503
$ws.Range("J6").Value = 0.4291715810931503
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.003554
$ws.Range("N6").Value = 3.010662
$ws.Range("O6").Value = 0.475327031240749
$ws.Range("P6").Value = 0.4753270312407492
$ws.Range("Q6").Value = 0.105288871464
$ws.Range("R6").Value = 0.9475998431760001
$ws.Range("S6").Value = 0.2039968535339055
$ws.Range("T6").Value = 0.2039968535339056

# Row 7: MuSCs -> Resolving-Mac
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ccl28"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.104916
$ws.Range("H7").Value = 0.314748
$ws.Range("I7").Value = 0.4291715810931503
$ws.Range("J7").Value = 0.4291715810931503
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.048954
$ws.Range("N7").Value = 0.146862
$ws.Range("O7").Value = 0.02318675376448066
$ws.Range("P7").Value = 0.02318675376448067
$ws.Range("Q7").Value = 0.005136057864
$ws.Range("R7").Value = 0.046224520776
$ws.Range("S7").Value = 0.009951095773519721
$ws.Range("T7").Value = 0.009951095773519724

